# Replace each of the 100 arithmetic-expression cells in the single
# 20x5 table with its updated value, in document (row-major) order.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    "61-59=",
    "55+38=",
    "6+76=",
    "38+59=",
    "41-19=",
    "85-26=",
    "51-39=",
    "57-38=",
    "19+2=",
    "10-8=",
    "43-25=",
    "17+58=",
    "86-49=",
    "32-26=",
    "34-7=",
    "47+45=",
    "95-28=",
    "65-17=",
    "59+12=",
    "22+49=",
    "74-45=",
    "58+19=",
    "55-17=",
    "83-19=",
    "28+44=",
    "17+44=",
    "5+16=",
    "66-38=",
    "75-57=",
    "53-7=",
    "74-57=",
    "44+8=",
    "16+69=",
    "38+14=",
    "60-54=",
    "33+19=",
    "64-57=",
    "74-57=",
    "67-38=",
    "37+36=",
    "54-17=",
    "59+29=",
    "12+49=",
    "24+27=",
    "53-48=",
    "70-57=",
    "73-68=",
    "92-84=",
    "71-47=",
    "20-1=",
    "3+78=",
    "42+29=",
    "38+35=",
    "20-7=",
    "75-19=",
    "17+45=",
    "29+12=",
    "83-78=",
    "28+36=",
    "25+39=",
    "39+12=",
    "82-8=",
    "45+16=",
    "19+39=",
    "58+16=",
    "80-1=",
    "60-57=",
    "57+5=",
    "28+14=",
    "70-37=",
    "24+58=",
    "45-6=",
    "43-34=",
    "18+79=",
    "35+58=",
    "51-42=",
    "50-38=",
    "70-61=",
    "29+32=",
    "94-48=",
    "43-27=",
    "52-18=",
    "82-3=",
    "90-49=",
    "83-18=",
    "76+5=",
    "34-8=",
    "47+47=",
    "91-44=",
    "83-15=",
    "38-19=",
    "84-75=",
    "56-28=",
    "54+19=",
    "29+43=",
    "19+23=",
    "32-27=",
    "64-28=",
    "95-26=",
    "93-65="
)

$i = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        $t.Cell($r, $c).Range.Text = $newValues[$i]
        $i = $i + 1
    }
}
